$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "fhkxk898"
$ws.Range("B2").Value = 231027163
$ws.Range("C2").Value = "zexfwlp69"
$ws.Range("D2").Value = 'ma&3$EC8'
$ws.Range("F2").Value = "tUZgzeYb"
$ws.Range("G2").Value = "ZeiM"
